# ---------------------------------------------------------------------------
# "Testing + model polishing" — updated syngas production and upstream
# emissions; new testing sheet with HIsarna energy mix; disconnect coke
# from pellet production (on the "removals" sheet).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$emissions = $wb.Worksheets.Item("emissions")

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet between "emissions" and "removals". Excel
#    auto-names it "Sheet1" and assigns the next free sheetId (3).
#    NOTE: worksheet references resolve by tab *position*, so any sheet
#    whose position shifts because of this insert (i.e. "removals", which
#    moves from slot 2 to slot 3) must be re-fetched afterwards.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $emissions)
$emissions = $wb.Worksheets.Item("emissions")
$removals  = $wb.Worksheets.Item("removals")

# ===========================================================================
# 2. "emissions" sheet — small data tweaks.
# ===========================================================================

# New value for "wood oven dry - IPCC" row (was blank).
$emissions.Range("B8").Value = 0

# "electricity PROXY - CN 2016" emissions factor revised down slightly.
$emissions.Range("B14").Value = 0.19400000000000001

# "wood (wet)" renamed / re-quantified -> "wood (20% moisture)", and its
# font colour is explicitly set to black (creates a new font/cellXf).
$emissions.Range("A15").Value = "wood (20% moisture)"
$emissions.Range("A15").Font.Color = 0
$emissions.Range("B15").Value = 0.6

# ===========================================================================
# 3. New "Sheet1" — HIsarna energy-mix testing table.
# ===========================================================================

$newSheet.Range("A1").Value = "substance"
$newSheet.Range("B1").Value = "CO2"

$newSheet.Range("A2").Value = "charcoal - IPCC"
$newSheet.Range("B2").Value = 2.992
$newSheet.Range("C2").Value = 4.54

$newSheet.Range("A3").Value = "coal bituminous - IPCC"
$newSheet.Range("B3").Value = 0.24137999999999998
$newSheet.Range("C3").Value = [double]"2.06E-2"

$newSheet.Range("A4").Value = "coal coking - IPCC"
$newSheet.Range("B4").Value = 0.24137999999999998
$newSheet.Range("C4").Value = [double]"2.06E-2"

$newSheet.Range("A5").Value = "natural gas - IPCC"
$newSheet.Range("B5").Value = 0.35599999999999998
$newSheet.Range("C5").Value = 0

$newSheet.Range("A6").Value = "wood oven dry - IPCC"
$newSheet.Range("B6").Value = 0
$newSheet.Range("C6").Value = 0

$newSheet.Range("A7").Value = "wood air dry - IPCC"
$newSheet.Range("B7").Value = 0.36520000000000002
$newSheet.Range("C7").Value = 0

$newSheet.Range("A8").Value = "electricity PROXY - EU 2016"
$newSheet.Range("B8").Value = [double]"4.8364000000000002E-3"
$newSheet.Range("C8").Value = [double]"8.8700000000000001E-5"

$newSheet.Range("A9").Value = "iron ore"
$newSheet.Range("B9").Value = [double]"6.2744500000000009E-2"
$newSheet.Range("C9").Formula = "=(0.0114+0.123)/2"

$newSheet.Range("A10").Value = "electricity PROXY - decarbonized"
$newSheet.Range("B10").Value = 0
$newSheet.Range("C10").Value = [char]8211

$newSheet.Range("A11").Value = "electricity PROXY - EU 2016"
$newSheet.Range("B11").Value = [double]"7.2999999999999995E-2"
$newSheet.Range("C11").Value = [char]8211

$newSheet.Range("A12").Value = "electricity PROXY - CN 2016"
$newSheet.Range("B12").Value = 0.19400000000000001
$newSheet.Range("C12").Value = [char]8211

$newSheet.Range("A13").Value = "wood (wet)"
$newSheet.Range("C13").Value = 1.2

$newSheet.Range("A14").Value = "solvent (MEA)"
$newSheet.Range("B14").Value = 1.8638999999999999
$newSheet.Range("C14").Value = [char]8211

# ===========================================================================
# 4. "removals" sheet — drop the hidden scratch rows (charcoal-2050 and the
#    never-used electricity proxy rows), and re-quantify the wood row.
# ===========================================================================

$removals.Range("A1:B24").EntireRow.Delete()

$removals.Range("A1").Value = "substance"
$removals.Range("B1").Value = "CO2 removed"

$removals.Range("A2").Value = "meta-notes"

$removals.Range("A3").Value = "meta-units"
$removals.Range("B3").Value = "t CO2 / t fuel"

$removals.Range("A4").Value = "charcoal - IPCC"
$removals.Range("B4").Value = 4.54

$removals.Range("A5").Value = "coal bituminous - IPCC"
$removals.Range("B5").Value = [double]"2.06E-2"

$removals.Range("A6").Value = "coal coking - IPCC"
$removals.Range("B6").Value = [double]"2.06E-2"

$removals.Range("A7").Value = "natural gas - IPCC"
$removals.Range("B7").Value = 0

$removals.Range("A8").Value = "wood oven dry - IPCC"
$removals.Range("B8").Value = 0

$removals.Range("A9").Value = "wood air dry - IPCC"
$removals.Range("B9").Value = 0

$removals.Range("A10").Value = "CaCO3"
$removals.Range("B10").Value = [double]"8.8700000000000001E-5"

$removals.Range("A11").Value = "iron ore"
$removals.Range("B11").Formula = "=(0.0114+0.123)/2"

$removals.Range("A12").Value = "wood (20% moisture)"
$removals.Range("B12").Value = 3.54

# ===========================================================================
# 5. Styles / cosmetics to match: row heights + column widths + header
#    styles on the new sheet, active-cell selections per sheet.
# ===========================================================================

$newSheet.Range("A1:B1").Font.Bold = $false
$newSheet.Rows.Item(1).RowHeight = 16

$newSheet.Activate()
$newSheet.Range("B13").Select() | Out-Null

$removals.Activate()
$removals.Range("A12").Select() | Out-Null

$emissions.Activate()
$emissions.Range("A22").Select() | Out-Null

Write-Output "edit complete"
